$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 616238.5361209477
$ws.Range("G2").Value = 616266.4689598368

# Row 3
$ws.Range("B3").Value = 0.6753301551942219
$ws.Range("C3").Value = 114.8270160096505
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 616238.5361209477
$ws.Range("G3").Value = 616354.8439567491

# Row 4
$ws.Range("B4").Value = 0.0001488876196638067
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 616238.5361209477
$ws.Range("G4").Value = 616241.0095540552

$wb.Save()
